# New crime data collected — weekly update for the 122nd Precinct CompStat
# report: bump the report volume/number and the covered week's dates in the
# header, then refresh the crime-complaint figures (week-to-date, 28-day,
# year-to-date, 2-year comparisons and their % changes) for rows 16-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: "Volume 31   Number  3" -> "...Number  4" and the covered week
# "1/15/2024 .. 1/21/2024" -> "1/22/2024 .. 1/28/2024"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/22/2024  Through  1/28/2024"

# ---------------------------------------------------------------------
# Helper: set a cell's value, then clone the number-format/alignment of a
# same-column "template" cell that already carries the style we need, so
# the destination ends up sharing the existing style (instead of Excel
# minting a brand-new near-duplicate xf record).
# ---------------------------------------------------------------------
function Set-CellWithStyle($target, $value, $templateAddr) {
    $target.Value = $value
    $ws.Range($templateAddr).Copy() | Out-Null
    $target.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Cells that change TYPE/STYLE (numeric <-> blank-placeholder text), so
# they need the value+format treatment above.
# ---------------------------------------------------------------------
Set-CellWithStyle $ws.Range("C16") 1        "C17"
Set-CellWithStyle $ws.Range("C18") "'0"     "C22"
Set-CellWithStyle $ws.Range("D18") "'0"     "D22"
Set-CellWithStyle $ws.Range("E18") "***.*"  "E22"
Set-CellWithStyle $ws.Range("F23") "'0"     "F22"
Set-CellWithStyle $ws.Range("M23") -100     "M16"
Set-CellWithStyle $ws.Range("C26") 1        "C17"
Set-CellWithStyle $ws.Range("F26") 1        "F16"
Set-CellWithStyle $ws.Range("I26") 1        "I16"
Set-CellWithStyle $ws.Range("D27") "'0"     "D22"
Set-CellWithStyle $ws.Range("E27") "***.*"  "E22"
Set-CellWithStyle $ws.Range("F30") "'0"     "F22"

# ---------------------------------------------------------------------
# Remaining cells only change value (style/number-format unchanged).
# ---------------------------------------------------------------------

# Row 16
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = -50

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 15
$ws.Range("K17").Value = -33.333333333333
$ws.Range("L17").Value = 42.857142857142
$ws.Range("M17").Value = -16.666666666666
$ws.Range("N17").Value = -44.444444444444

# Row 18
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("M18").Value = -71.428571428571
$ws.Range("N18").Value = -95.348837209302

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("E19").Value = 77.777777777777
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -2.127659574468
$ws.Range("I19").Value = 46
$ws.Range("J19").Value = 47
$ws.Range("K19").Value = -2.127659574468
$ws.Range("L19").Value = 53.333333333333
$ws.Range("M19").Value = 70.370370370370
$ws.Range("N19").Value = -37.837837837837

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = -60
$ws.Range("L20").Value = -66.666666666666
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -98.4

# Row 21
$ws.Range("C21").Value = 24
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -18.072289156626
$ws.Range("I21").Value = 68
$ws.Range("J21").Value = 83
$ws.Range("K21").Value = -18.072289156626
$ws.Range("L21").Value = 25.925925925925
$ws.Range("M21").Value = -1.449275362318
$ws.Range("N21").Value = -85.152838427947

# Row 23
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = -100

# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -36.363636363636
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = -9.090909090909
$ws.Range("I24").Value = 70
$ws.Range("J24").Value = 77
$ws.Range("K24").Value = -9.090909090909
$ws.Range("L24").Value = 22.807017543859
$ws.Range("M24").Value = -39.655172413793

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -29.411764705882
$ws.Range("I25").Value = 12
$ws.Range("J25").Value = 17
$ws.Range("K25").Value = -29.411764705882
$ws.Range("L25").Value = -20
$ws.Range("M25").Value = -50

# Row 26
$ws.Range("L26").Value = -50

# Row 27
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -75
